$wb = $excel.ActiveWorkbook

# --- mobility sheet ---
$mobility = $wb.Worksheets.Item("mobility")
$mobility.Range("B24").Value = 420000
$mobility.Range("B26").HorizontalAlignment = -4108   # xlCenter -- promotes B26's style (font2) to a new centered xf
$mobility.Activate() | Out-Null
$mobility.Range("B26").Select() | Out-Null

# --- infrastructure sheet ---
$infrastructure = $wb.Worksheets.Item("infrastructure")
$infrastructure.Range("B16").Value = 156250.76923076919
$infrastructure.Range("B17").Value = 701400
$infrastructure.Activate() | Out-Null
$infrastructure.Range("B17").Select() | Out-Null

# --- categories sheet ---
$categories = $wb.Worksheets.Item("categories")
$categories.Range("B2:B6").NumberFormat = "0.0%"
$categories.Range("D2:D6").NumberFormat = "0%"
$categories.Activate() | Out-Null
$categories.Range("B6").Select() | Out-Null
